# Swap the first two comma-separated names in the "Recorded By" (column G)
# values on the active sheet, for every data row that has two or more
# comma-separated entries. Entries with only a single value are left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = 7
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -isnot [string]) { continue }

    $parts = $val -split ', '
    if ($parts.Count -ge 2) {
        $tmp = $parts[0]
        $parts[0] = $parts[1]
        $parts[1] = $tmp
        $cell.Value = [string]::Join(', ', $parts)
    }
}
